$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.020.67'
$ws.Range('E2').Value = '  +1.06%  '

$ws.Range('D3').Value = '2.693.43'
$ws.Range('E3').Value = '  +1.71%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '611.15'
$ws.Range('E5').Value = '  +1.22%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.44'
$ws.Range('E6').Value = '  +0.94%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('E9').Value = '  +3.99%  '

$ws.Range('E10').Value = '  +3.56%  '

$ws.Range('E11').Value = '  -1.46%  '

$ws.Range('E12').Value = '  +0.26%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000208'
$ws.Range('E13').Value = '  +9.71%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.09'
$ws.Range('E14').Value = '  +2.26%  '

$ws.Range('D15').Value = '3.177.64'
$ws.Range('E15').Value = '  +1.74%  '

$ws.Range('D16').Value = '65.865.48'
$ws.Range('E16').Value = '  +1.13%  '

$ws.Range('D17').Value = '2.683.17'
$ws.Range('E17').Value = '  +1.07%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.79'
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.89'
$ws.Range('E19').Value = '  -0.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.74'
$ws.Range('E20').Value = '  +4.79%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '358.10'
$ws.Range('E21').Value = '  -0.59%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.34'
$ws.Range('E22').Value = '  +2.76%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('E24').Value = '  +16.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.92'
$ws.Range('E25').Value = '  +4.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.63'
$ws.Range('E26').Value = '  -5.39%  '

$ws.Range('E27').Value = '  +0.22%  '

$ws.Range('E28').Value = '  +3.56%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.29'
$ws.Range('E29').Value = '  -0.54%  '

$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '532.83'
$ws.Range('E32').Value = '  -3.46%  '

$ws.Range('E33').Value = '  -1.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.67'
$ws.Range('E34').Value = '  +3.92%  '

$ws.Range('E35').Value = '  -2.51%  '

$ws.Range('E36').Value = '  +0.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.75'
$ws.Range('E37').Value = '  +0.84%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.54'
$ws.Range('E38').Value = '  +0.16%  '

$ws.Range('E39').Value = '  -1.40%  '

$ws.Range('E40').Value = '  +0.07%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.96'
$ws.Range('E42').Value = '  +0.64%  '

$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.43'
$ws.Range('E43').Value = '  -0.79%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.16'
$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0637'
$ws.Range('E45').Value = '  +2.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.82'
$ws.Range('E46').Value = '  +1.97%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.30'
$ws.Range('E47').Value = '  +0.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0267'
$ws.Range('E48').Value = '  +1.06%  '

$ws.Range('E49').Value = '  +0.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.90'
$ws.Range('E50').Value = '  +5.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0997'
$ws.Range('E51').Value = '  +1.46%  '
